$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "317.65"
    "E2" = "2.80%"
    "D3" = "41.40"
    "E3" = "1.53%"
    "D4" = "5.249"
    "E4" = "2.36%"
    "D5" = "0.07739"
    "E5" = "1.44%"
    "D6" = "1.705"
    "E6" = "5.16%"
    "D7" = "0.9517"
    "E7" = "4.63%"
    "E8" = "-2.94%"
    "D9" = "0.1263"
    "E9" = "5.65%"
    "D10" = "0.1840"
    "E10" = "1.88%"
    "D11" = "0.09249"
    "E11" = "1.21%"
    "D12" = "0.04390"
    "E12" = "3.16%"
    "D13" = "0.1051"
    "E13" = "0.60%"
    "D14" = "0.001283"
    "E14" = "2.07%"
    "D15" = "0.005959"
    "E15" = "2.91%"
    "E16" = "-0.05%"
    "D17" = "4.322"
    "E17" = "0.97%"
    "D18" = "0.3352"
    "E18" = "2.90%"
    "D19" = "7.681"
    "E19" = "11.20%"
    "D20" = "0.1351"
    "E20" = "-4.28%"
    "D21" = "0.2820"
    "E21" = "4.20%"
    "D22" = "0.04015"
    "E22" = "-0.82%"
    "D23" = "0.001265"
    "E23" = "-0.57%"
    "D24" = "0.004124"
    "E24" = "0.34%"
    "D25" = "0.0001269"
    "E25" = "-0.23%"
    "D38" = "0.02545"
    "E38" = "5.10%"
    "D39" = "0.05358"
    "E39" = "2.62%"
    "D40" = "0.007781"
    "E40" = "-0.21%"
    "E41" = "1.54%"
    "D42" = "0.007309"
    "E42" = "7.36%"
    "D43" = "0.001977"
    "E43" = "2.28%"
    "D44" = "0.007560"
    "E44" = "-6.53%"
    "D45" = "0.3431"
    "E45" = "11.65%"
    "D46" = "0.00006679"
    "E46" = "-3.22%"
    "E47" = "-0.27%"
    "D48" = "0.2177"
    "E48" = "119.41%"
    "D49" = "0.004198"
    "E49" = "39.81%"
    "D50" = "0.00002099"
    "E50" = "-0.27%"
    "D51" = "0.0001999"
    "E51" = "-0.27%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}
